# Update RAD Phase 3 Test Cases and test data for Estate Tax.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Date column (B) for the existing Personal Income Tax rows
# with the new test-run timestamps.
$ws.Range("B2").Value = "Fri Feb 02 19:27:22 EST 2024"
$ws.Range("B3").Value = "Fri Feb 02 19:27:35 EST 2024"
$ws.Range("B5").Value = "Fri Feb 02 19:27:47 EST 2024"

# Add the two new Estate Tax rows (Result / Date / Execute columns),
# matching the existing PaymentType/TaxType values already present in
# D6:E7.
$ws.Range("A6").Value = "Pass"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "Fri Feb 02 19:27:58 EST 2024"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "Y"

$ws.Range("A7").Value = "Pass"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "Fri Feb 02 19:28:10 EST 2024"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "Y"

# Update current selection to C5
$ws.Range("C5").Select()
